# New PO forecast model
# - Appends new weekly rows to "Weekly Quantity"
# - Appends a new monthly row to "Monthly Trend"
# - Adds a brand-new "PO Forecast" sheet with a Prophet-style ds/PO_Forecast series

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Weekly Quantity" — append two more weekly data points
# ---------------------------------------------------------------------------
$weekly = $wb.Worksheets.Item("Weekly Quantity")

$weekly.Cells.Item(3, 1).Value = 45662.99999999999
$weekly.Cells.Item(3, 1).NumberFormat = $weekly.Cells.Item(2, 1).NumberFormat
$weekly.Cells.Item(3, 2).Value = 14

$weekly.Cells.Item(4, 1).Value = 45676.99999999999
$weekly.Cells.Item(4, 1).NumberFormat = $weekly.Cells.Item(2, 1).NumberFormat
$weekly.Cells.Item(4, 2).Value = 6

# ---------------------------------------------------------------------------
# 2. "Monthly Trend" — append one more monthly data point
# ---------------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly Trend")

$monthly.Cells.Item(3, 1).Value = 45688.99999999999
$monthly.Cells.Item(3, 1).NumberFormat = $monthly.Cells.Item(2, 1).NumberFormat
$monthly.Cells.Item(3, 2).Value = 20

# ---------------------------------------------------------------------------
# 3. New "PO Forecast" sheet, placed after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

$forecast.Cells.Item(1, 1).Value = "ds"
$forecast.Cells.Item(1, 2).Value = "PO_Forecast"

# match the bold/centered/bordered header styling used on the other sheets
$weekly.Range("A1:B1").Copy()
$forecast.Range("A1:B1").PasteSpecial(-4122)

$forecastData = @(
    @(45641.99999999999, 10),
    @(45662.99999999999, 9),
    @(45676.99999999999, 9),
    @(45683.99999999999, 9),
    @(45690.99999999999, 8),
    @(45697.99999999999, 8),
    @(45704.99999999999, 8),
    @(45711.99999999999, 8),
    @(45718.99999999999, 8),
    @(45725.99999999999, 7),
    @(45732.99999999999, 7)
)

$row = 2
foreach ($point in $forecastData) {
    $forecast.Cells.Item($row, 1).Value = $point[0]
    $forecast.Cells.Item($row, 1).NumberFormat = $weekly.Cells.Item(2, 1).NumberFormat
    $forecast.Cells.Item($row, 2).Value = $point[1]
    $row++
}
